$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038048763168817
$ws.Range("D2").Value = 1.041673378081424
$ws.Range("E2").Value = 1.051311232318664
$ws.Range("F2").Value = 1.058170631353529
$ws.Range("I2").Value = 1.040545273721847
$ws.Range("J2").Value = 1.043148315267115
$ws.Range("K2").Value = 1.044452227391167
$ws.Range("L2").Value = 1.054063086190362
$ws.Range("M2").Value = 1.060903596684663
$ws.Range("N2").Value = 1.018303447963193
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038946480183102
$ws.Range("D3").Value = 1.04236117998214
$ws.Range("E3").Value = 1.052229482881794
$ws.Range("F3").Value = 1.059230317303665
$ws.Range("I3").Value = 1.040785481165968
$ws.Range("J3").Value = 1.043690924032695
$ws.Range("K3").Value = 1.044951027510634
$ws.Range("L3").Value = 1.054793685917362
$ws.Range("M3").Value = 1.061776639859837
$ws.Range("N3").Value = 1.018485015655082
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039527704346274
$ws.Range("D4").Value = 1.042806386335034
$ws.Range("E4").Value = 1.052824405708282
$ws.Range("F4").Value = 1.059917077632485
$ws.Range("I4").Value = 1.040939679152881
$ws.Range("J4").Value = 1.044041727527264
$ws.Range("K4").Value = 1.045273259368733
$ws.Range("L4").Value = 1.055266554768278
$ws.Range("M4").Value = 1.062342026063751
$ws.Range("N4").Value = 1.018602357037038
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039772131637591
$ws.Range("D5").Value = 1.042993586251797
$ws.Range("E5").Value = 1.053074690370467
$ws.Range("F5").Value = 1.060206047322036
$ws.Range("I5").Value = 1.041004208380006
$ws.Range("J5").Value = 1.044189132443138
$ws.Range("K5").Value = 1.045408598958205
$ws.Range("L5").Value = 1.055465376993169
$ws.Range("M5").Value = 1.062579825601461
$ws.Range("N5").Value = 1.018651652276042
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039813176754797
$ws.Range("D6").Value = 1.043025019964411
$ws.Range("E6").Value = 1.053116724733461
$ws.Range("F6").Value = 1.060254581562385
$ws.Range("I6").Value = 1.041015025769446
$ws.Range("J6").Value = 1.044213878077817
$ws.Range("K6").Value = 1.045431315614216
$ws.Range("L6").Value = 1.055498761740685
$ws.Range("M6").Value = 1.062619759678506
$ws.Range("N6").Value = 1.018659927096378
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039530970080082
$ws.Range("D7").Value = 1.042808887573001
$ws.Range("E7").Value = 1.05282774932197
$ws.Range("F7").Value = 1.059920937856532
$ws.Range("I7").Value = 1.040940542557842
$ws.Range("J7").Value = 1.044043697446425
$ws.Range("K7").Value = 1.045275068281929
$ws.Range("L7").Value = 1.055269211331228
$ws.Range("M7").Value = 1.062345203117302
$ws.Range("N7").Value = 1.018603015860707
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038352079559332
$ws.Range("D8").Value = 1.041905791560165
$ws.Range("E8").Value = 1.051621402672964
$ws.Range("F8").Value = 1.058528534388308
$ws.Range("I8").Value = 1.040626707749272
$ws.Range("J8").Value = 1.043331754155518
$ws.Range("K8").Value = 1.044620907392605
$ws.Range("L8").Value = 1.05430996990713
$ws.Range("M8").Value = 1.061198547978797
$ws.Range("N8").Value = 1.018364839438335
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036277382993428
$ws.Range("D9").Value = 1.040315648741054
$ws.Range("E9").Value = 1.049501489849605
$ws.Range("F9").Value = 1.056083207628905
$ws.Range("I9").Value = 1.040064279027695
$ws.Range("J9").Value = 1.042074955876281
$ws.Range("K9").Value = 1.043464212290454
$ws.Range("L9").Value = 1.052620644057561
$ws.Range("M9").Value = 1.059181633752955
$ws.Range("N9").Value = 1.017944046688513
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034896100306422
$ws.Range("D10").Value = 1.039256458628659
$ws.Range("E10").Value = 1.048092203118264
$ws.Range("F10").Value = 1.054458614596557
$ws.Range("I10").Value = 1.039683033498416
$ws.Range("J10").Value = 1.041235621167115
$ws.Range("K10").Value = 1.042690462465924
$ws.Range("L10").Value = 1.051495146810815
$ws.Range("M10").Value = 1.057839533630809
$ws.Range("N10").Value = 1.017662802325523
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034298440416679
$ws.Range("D11").Value = 1.038798048889095
$ws.Range("E11").Value = 1.047482927431447
$ws.Range("F11").Value = 1.053756495767384
$ws.Range("I11").Value = 1.039516464526057
$ws.Range("J11").Value = 1.040871842329798
$ws.Range("K11").Value = 1.042354810511364
$ws.Range("L11").Value = 1.051007977423573
$ws.Range("M11").Value = 1.057258997572115
$ws.Range("N11").Value = 1.017540854902196
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0340765106684
$ws.Range("D12").Value = 1.038627810454006
$ws.Range("E12").Value = 1.047256759930871
$ws.Range("F12").Value = 1.053495900057752
$ws.Range("I12").Value = 1.039454370509082
$ws.Range("J12").Value = 1.040736668223326
$ws.Range("K12").Value = 1.042230043211273
$ws.Range("L12").Value = 1.050827048791507
$ws.Range("M12").Value = 1.057043452096203
$ws.Range("N12").Value = 1.017495533479965
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034124112271922
$ws.Range("D13").Value = 1.038664325555356
$ws.Range("E13").Value = 1.047305267064366
$ws.Range("F13").Value = 1.053551789542066
$ws.Range("I13").Value = 1.039467699963763
$ws.Range("J13").Value = 1.040765665797486
$ws.Range("K13").Value = 1.042256810334033
$ws.Range("L13").Value = 1.050865857324959
$ws.Range("M13").Value = 1.057089683177443
$ws.Range("N13").Value = 1.017505256192127
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034280094246984
$ws.Range("D14").Value = 1.038783976198164
$ws.Range("E14").Value = 1.047464229398341
$ws.Range("F14").Value = 1.053734950707922
$ws.Range("I14").Value = 1.03951133635887
$ws.Range("J14").Value = 1.040860669820283
$ws.Range("K14").Value = 1.042344499063036
$ws.Range("L14").Value = 1.050993021234268
$ws.Range("M14").Value = 1.057241178644198
$ws.Range("N14").Value = 1.017537109121634
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034376208895162
$ws.Range("D15").Value = 1.038857701632032
$ws.Range("E15").Value = 1.047562190537114
$ws.Range("F15").Value = 1.053847829210468
$ws.Range("I15").Value = 1.039538192664013
$ws.Range("J15").Value = 1.040919198274355
$ws.Range("K15").Value = 1.04239851492266
$ws.Range("L15").Value = 1.051071374824668
$ws.Range("M15").Value = 1.057334532151121
$ws.Range("N15").Value = 1.017556731497202
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034935774473364
$ws.Range("D16").Value = 1.039286886659033
$ws.Range("E16").Value = 1.04813265895506
$ws.Range("F16").Value = 1.054505240248093
$ws.Range("I16").Value = 1.039694056846389
$ws.Range("J16").Value = 1.0412597568497
$ws.Range("K16").Value = 1.042712725738914
$ws.Range("L16").Value = 1.051527482479338
$ws.Range("M16").Value = 1.057878074668842
$ws.Range("N16").Value = 1.017670892089493
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035286894882906
$ws.Range("D17").Value = 1.039556164674206
$ws.Range("E17").Value = 1.048490755169048
$ws.Range("F17").Value = 1.054917976510151
$ws.Range("I17").Value = 1.039791428562132
$ws.Range("J17").Value = 1.041473289443452
$ws.Range("K17").Value = 1.042909658355952
$ws.Range("L17").Value = 1.051813635197283
$ws.Range("M17").Value = 1.058219186446698
$ws.Range("N17").Value = 1.017742457602872
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035491740313971
$ws.Range("D18").Value = 1.039713251757237
$ws.Range("E18").Value = 1.048699718674547
$ws.Range("F18").Value = 1.055158847872948
$ws.Range("I18").Value = 1.039848080312408
$ws.Range("J18").Value = 1.041597806399254
$ws.Range("K18").Value = 1.043024466596588
$ws.Range("L18").Value = 1.051980560378838
$ws.Range("M18").Value = 1.058418209207269
$ws.Range("N18").Value = 1.017784184443884
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035561594549423
$ws.Range("D19").Value = 1.039766818046082
$ws.Range("E19").Value = 1.048770985406103
$ws.Range("F19").Value = 1.05524100067949
$ws.Range("I19").Value = 1.039867372734589
$ws.Range("J19").Value = 1.04164025784164
$ws.Range("K19").Value = 1.043063603167002
$ws.Range("L19").Value = 1.05203748043005
$ws.Range("M19").Value = 1.05848608063827
$ws.Range("N19").Value = 1.017798409475288
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035249218548273
$ws.Range("D20").Value = 1.039527271430051
$ws.Range("E20").Value = 1.048452325305388
$ws.Range("F20").Value = 1.054873680435145
$ws.Range("I20").Value = 1.039780996339993
$ws.Range("J20").Value = 1.041450382823588
$ws.Range("K20").Value = 1.042888535469068
$ws.Range("L20").Value = 1.051782931938573
$ws.Range("M20").Value = 1.058182582371449
$ws.Range("N20").Value = 1.017734780964711
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034234159565359
$ws.Range("D21").Value = 1.038748741079618
$ws.Range("E21").Value = 1.047417414951426
$ws.Range("F21").Value = 1.053681008711598
$ws.Range("I21").Value = 1.039498492676451
$ws.Range("J21").Value = 1.040832694889389
$ws.Range("K21").Value = 1.042318679428771
$ws.Range("L21").Value = 1.05095557385804
$ws.Range("M21").Value = 1.05719656447226
$ws.Range("N21").Value = 1.017527729907505
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033596344120913
$ws.Range("D22").Value = 1.038259453338478
$ws.Range("E22").Value = 1.046767563031548
$ws.Range("F22").Value = 1.052932300709601
$ws.Range("I22").Value = 1.039319582141273
$ws.Range("J22").Value = 1.040444038502235
$ws.Range("K22").Value = 1.041959860999659
$ws.Range("L22").Value = 1.050435542145304
$ws.Range("M22").Value = 1.056577145566912
$ws.Range("N22").Value = 1.017397405835978
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03393442464188
$ws.Range("D23").Value = 1.038518814139866
$ws.Range("E23").Value = 1.047111982071306
$ws.Range("F23").Value = 1.053329093534535
$ws.Range("I23").Value = 1.039414548027436
$ws.Range("J23").Value = 1.040650099954316
$ws.Range("K23").Value = 1.042150127120031
$ws.Range("L23").Value = 1.050711205194337
$ws.Range("M23").Value = 1.056905460706451
$ws.Range("N23").Value = 1.01746650652148
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035266242729837
$ws.Range("D24").Value = 1.039540326976252
$ws.Range("E24").Value = 1.048469689822286
$ws.Range("F24").Value = 1.054893695525343
$ws.Range("I24").Value = 1.039785710655633
$ws.Range("J24").Value = 1.041460733441113
$ws.Range("K24").Value = 1.042898080175495
$ws.Range("L24").Value = 1.051796805366687
$ws.Range("M24").Value = 1.058199121998211
$ws.Range("N24").Value = 1.017738249756873
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03681342044104
$ws.Range("D25").Value = 1.040726585069395
$ws.Range("E25").Value = 1.050048840458469
$ws.Range("F25").Value = 1.05671439673117
$ws.Range("I25").Value = 1.040210792082287
$ws.Range("J25").Value = 1.042400131157017
$ws.Range("K25").Value = 1.04376371156482
$ws.Range("L25").Value = 1.053057252853841
$ws.Range("M25").Value = 1.05970261675221
$ws.Range("N25").Value = 1.018052959389083
